# Auto commit at 2026-01-19  9:48:32.37
# Updates the monthly "Metrics" figures (column B, rows 2-13) and lets the
# dependent formulas on the "today" sheet recalculate. Also restores the
# per-sheet selection / active-tab state recorded at save time.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the raw metric values on the "Metrics" sheet.
#    (Sheet "today" references these via formulas such as "=Metrics!B2",
#    so its cached values refresh automatically once these change.)
# ---------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 332199
$wsMetrics.Range("B3").Value  = 246828.34000000003
$wsMetrics.Range("B4").Value  = 84335.16
$wsMetrics.Range("B5").Value  = 13616
$wsMetrics.Range("B6").Value  = 5968069.7300000004
$wsMetrics.Range("B7").Value  = 5017545.97
$wsMetrics.Range("B8").Value  = 1748426.98
$wsMetrics.Range("B9").Value  = 233893
$wsMetrics.Range("B10").Value = 34433450.719999999
$wsMetrics.Range("B11").Value = 32292821.129999999
$wsMetrics.Range("B12").Value = 12030149.02
$wsMetrics.Range("B13").Value = 1331523

# ---------------------------------------------------------------------
# 2) Restore the selection recorded on the "today" sheet, then switch
#    back to "Metrics" so it ends up as the active/selected tab (matching
#    the saved workbook state).
# ---------------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Activate() | Out-Null
$wsToday.Range("C31").Select() | Out-Null

$wsMetrics.Activate() | Out-Null
$wsMetrics.Range("D8").Select() | Out-Null
